$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, taken from the commit diff
# (updated Price (D) and Volume(1h) (E) columns of the cryptos table).
$updates = @{
    'D2' = '55.550.35'
    'E2' = '  +3.66%  '
    'D3' = '2.499.97'
    'E3' = '  +7.03%  '
    'E4' = '  +0.17%  '
    'D5' = '482.82'
    'E5' = '  +9.40%  '
    'D6' = '140.95'
    'E6' = '  +14.90%  '
    'D7' = '1.00'
    'E7' = '  +0.47%  '
    'D8' = '0.511'
    'E8' = '  +7.42%  '
    'D9' = '2.495.93'
    'E9' = '  +6.79%  '
    'D10' = '0.0988'
    'E10' = '  +7.20%  '
    'D11' = '5.49'
    'E11' = '  +3.78%  '
    'E12' = '  +6.04%  '
    'E13' = '  +0.35%  '
    'D14' = '2.931.41'
    'E14' = '  +7.51%  '
    'D15' = '55.593.96'
    'E15' = '  +3.76%  '
    'D16' = '20.61'
    'E16' = '  +7.12%  '
    'E17' = '  +14.25%  '
    'D18' = '2.495.99'
    'E18' = '  +7.19%  '
    'D19' = '4.38'
    'E19' = '  +9.36%  '
    'D20' = '320.24'
    'E20' = '  +5.59%  '
    'D21' = '10.04'
    'E21' = '  +8.76%  '
    'D22' = '1.00'
    'E22' = '  +0.24%  '
    'D23' = '5.70'
    'E23' = '  +5.32%  '
    'D24' = '57.97'
    'E24' = '  +4.15%  '
    'D25' = '0.167'
    'E25' = '  +8.55%  '
    'D26' = '0.410'
    'E26' = '  +10.58%  '
    'E27' = '  +0.90%  '
    'D28' = '2.606.78'
    'E28' = '  +6.96%  '
    'D29' = '7.44'
    'E29' = '  +6.12%  '
    'D30' = '0.0₃0791'
    'E30' = '  +11.86%  '
    'E31' = '  +0.46%  '
    'D32' = '149.25'
    'E32' = '  +2.77%  '
    'D33' = '18.16'
    'E33' = '  +5.21%  '
    'E34' = '  +9.45%  '
    'D35' = '5.18'
    'E35' = '  +9.78%  '
    'E36' = '  +4.00%  '
    'E37' = '  +11.05%  '
    'D38' = '0.865'
    'E38' = '  +4.14%  '
    'D39' = '34.10'
    'E39' = '  +3.04%  '
    'D40' = '0.998'
    'E40' = '  +0.57%  '
    'D41' = '0.609'
    'E41' = '  +15.89%  '
    'E42' = '  +10.47%  '
    'D43' = '3.42'
    'E43' = '  +7.85%  '
    'E44' = '  +8.31%  '
    'D45' = '10.17'
    'E45' = '  -1.00%  '
    'D46' = '1.977.64'
    'E46' = '  +2.89%  '
    'D47' = '0.0906'
    'E47' = '  +8.82%  '
    'D48' = '4.63'
    'E48' = '  +16.79%  '
    'E49' = '  +6.74%  '
    'D50' = '252.73'
    'E50' = '  +33.40%  '
    'D51' = '17.59'
    'E51' = '  +12.41%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $originalStyle = $cell.Style
    # Force text interpretation so numeric-looking strings (e.g. "482.82")
    # are not silently coerced into numbers, matching the source data which
    # stores these as plain text (inline strings) in the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = $originalStyle
}
